$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values per repull of data
$ws.Range("F4").Value = -11
$ws.Range("F6").Value = 3
$ws.Range("F8").Value = -2
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -2
$ws.Range("F17").Value = -1
$ws.Range("F21").Value = 2
$ws.Range("F23").Value = -8
